$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 9,14
$arr[0,0] = 0.6504219999999999
$arr[0,1] = 1.951266
$arr[0,2] = 0.1521898546336546
$arr[0,3] = 0.1521898546336546
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 28.22405966666667
$arr[0,7] = 84.672179
$arr[0,8] = 0.3816548478108986
$arr[0,9] = 0.3816548478108986
$arr[0,10] = 18.35754933651267
$arr[0,11] = 165.217944028614
$arr[0,12] = 0.05808399580857024
$arr[0,13] = 0.05808399580857024
$arr[1,0] = 0.6504219999999999
$arr[1,1] = 1.951266
$arr[1,2] = 0.1521898546336546
$arr[1,3] = 0.1521898546336546
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 19.768727
$arr[1,7] = 59.306181
$arr[1,8] = 0.2673191094302723
$arr[1,9] = 0.2673191094302723
$arr[1,10] = 12.858014952794
$arr[1,11] = 115.722134575146
$arr[1,12] = 0.04068325640499115
$arr[1,13] = 0.04068325640499115
$arr[2,0] = 0.6504219999999999
$arr[2,1] = 1.951266
$arr[2,2] = 0.1521898546336546
$arr[2,3] = 0.1521898546336546
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 25.95900466666667
$arr[2,7] = 77.877014
$arr[2,8] = 0.351026042758829
$arr[2,9] = 0.351026042758829
$arr[2,10] = 16.88430773330267
$arr[2,11] = 151.958769599724
$arr[2,12] = 0.05342260242009322
$arr[2,13] = 0.05342260242009321
$arr[3,0] = 1.870396
$arr[3,1] = 5.611188
$arr[3,2] = 0.437647089654669
$arr[3,3] = 0.4376470896546689
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 28.22405966666667
$arr[3,7] = 84.672179
$arr[3,8] = 0.3816548478108986
$arr[3,9] = 0.3816548478108986
$arr[3,10] = 52.79016830429467
$arr[3,11] = 475.111514738652
$arr[3,12] = 0.1670301333970354
$arr[3,13] = 0.1670301333970354
$arr[4,0] = 1.870396
$arr[4,1] = 5.611188
$arr[4,2] = 0.437647089654669
$arr[4,3] = 0.4376470896546689
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 19.768727
$arr[4,7] = 59.306181
$arr[4,8] = 0.2673191094302723
$arr[4,9] = 0.2673191094302723
$arr[4,10] = 36.975347905892
$arr[4,11] = 332.778131153028
$arr[4,12] = 0.1169914302512366
$arr[4,13] = 0.1169914302512366
$arr[5,0] = 1.870396
$arr[5,1] = 5.611188
$arr[5,2] = 0.437647089654669
$arr[5,3] = 0.4376470896546689
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 25.95900466666667
$arr[5,7] = 77.877014
$arr[5,8] = 0.351026042758829
$arr[5,9] = 0.351026042758829
$arr[5,10] = 48.55361849251467
$arr[5,11] = 436.982566432632
$arr[5,12] = 0.1536255260063969
$arr[5,13] = 0.1536255260063969
$arr[6,0] = 1.752936
$arr[6,1] = 5.258808
$arr[6,2] = 0.4101630557116764
$arr[6,3] = 0.4101630557116764
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 28.22405966666667
$arr[6,7] = 84.672179
$arr[6,8] = 0.3816548478108986
$arr[6,9] = 0.3816548478108986
$arr[6,10] = 49.474970255848
$arr[6,11] = 445.274732302632
$arr[6,12] = 0.156540718605293
$arr[6,13] = 0.156540718605293
$arr[7,0] = 1.752936
$arr[7,1] = 5.258808
$arr[7,2] = 0.4101630557116764
$arr[7,3] = 0.4101630557116764
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 19.768727
$arr[7,7] = 59.306181
$arr[7,8] = 0.2673191094302723
$arr[7,9] = 0.2673191094302723
$arr[7,10] = 34.653313232472
$arr[7,11] = 311.879819092248
$arr[7,12] = 0.1096444227740445
$arr[7,13] = 0.1096444227740445
$arr[8,0] = 1.752936
$arr[8,1] = 5.258808
$arr[8,2] = 0.4101630557116764
$arr[8,3] = 0.4101630557116764
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 25.95900466666667
$arr[8,7] = 77.877014
$arr[8,8] = 0.351026042758829
$arr[8,9] = 0.351026042758829
$arr[8,10] = 45.504473804368
$arr[8,11] = 409.540264239312
$arr[8,12] = 0.1439779143323389
$arr[8,13] = 0.1439779143323389
$ws.Range("G2:T10").Value = $arr
Write-Output "done"
